$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure cells keep their original text representation (avoid Excel
# auto-converting numeric-looking strings like "239.21" into numbers,
# and preserve exact padding/precision such as "1.000" or "  -0.77%  ").

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.539.89"
$ws.Range("E2").Value = "  -0.77%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.912.27"
$ws.Range("E3").Value = "  -1.37%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "239.21"
$ws.Range("E5").Value = "  -1.28%  "
$ws.Range("E6").Value = "  -0.13%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4780"
$ws.Range("E7").Value = "  -2.13%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2836"
$ws.Range("E8").Value = "  -3.22%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06691"
$ws.Range("E9").Value = "  -2.78%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.67"
$ws.Range("E10").Value = "  -4.41%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "101.03"
$ws.Range("E11").Value = "  -4.16%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.917.55"
$ws.Range("E12").Value = "  -1.19%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07684"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.201"
$ws.Range("E14").Value = "  -2.14%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6675"
$ws.Range("E15").Value = "  -4.26%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "30.541.24"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "256.31"
$ws.Range("E17").Value = "  -6.94%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.001"
$ws.Range("E18").Value = "  -0.19%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007460"
$ws.Range("E19").Value = "  -3.24%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.63"
$ws.Range("E20").Value = "  -3.93%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.373"
$ws.Range("E21").Value = "  -1.19%  "
$ws.Range("E22").Value = "  +0.00%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.291"
$ws.Range("E23").Value = "  -2.63%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.309"
$ws.Range("E24").Value = "  -4.07%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "166.66"
$ws.Range("E25").Value = "  -0.68%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "19.13"
$ws.Range("E26").Value = "  -2.28%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.053"
$ws.Range("E27").Value = "  -4.99%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "4.730"
$ws.Range("E28").Value = "  +4.15%  "
$ws.Range("E29").Value = "  -0.76%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.1006"
$ws.Range("E30").Value = "  -3.28%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.509"
$ws.Range("E31").Value = "  -2.77%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.235"
$ws.Range("E32").Value = "  -2.83%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04711"
$ws.Range("E33").Value = "  -2.78%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7246"
$ws.Range("E34").Value = "  -3.18%  "
$ws.Range("E35").Value = "  -4.00%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.000"
$ws.Range("E36").Value = "  -0.09%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.708"
$ws.Range("E37").Value = "  -0.68%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01909"
$ws.Range("E38").Value = "  -3.94%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.611"
$ws.Range("E39").Value = "  -1.97%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.237"
$ws.Range("E40").Value = "  -3.04%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "74.54"
$ws.Range("E41").Value = "  -3.66%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.963"
$ws.Range("E42").Value = "  -6.16%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8601"
$ws.Range("E43").Value = "  -4.21%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "105.35"
$ws.Range("E44").Value = "  -2.59%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.001"
$ws.Range("E45").Value = "  +0.28%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4223"
$ws.Range("E46").Value = "  -4.16%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.337"
$ws.Range("E47").Value = "  -5.20%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.1197"
$ws.Range("E48").Value = "  -3.39%  "
$ws.Range("B49").Value = "Maker"
$ws.Range("C49").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "910.88"
$ws.Range("E49").Value = "  -8.03%  "
$ws.Range("B50").Value = "Elrond"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "34.66"
$ws.Range("E50").Value = "  -3.03%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.715"
$ws.Range("E51").Value = "  -4.81%  "
